# Generate Report for Handback
# Update the timestamp strings in the handback status report to reflect
# the latest xliff generation / handoff / handback times.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" on the Overview sheet, and the matching
# "Correspond Handoff Datetime" on the de-de sheet, both shared the same
# text (2016-08-22 10:44:27) and both move to 2016-08-22 10:45:25.
$wsOverview.Range("G2").Value = "2016-08-22 10:45:25"
$wsDeDe.Range("H2").Value     = "2016-08-22 10:45:25"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-22 10:45:20"
$wsZhCn.Range("K2").Value = "2016-08-22 10:45:38"

# de-de sheet: Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-08-22 10:45:45"
